# Reporte.xlsx - "Finalizacion de Proyecto" update
# - Refresh the report timestamp (title + last-row FECHA)
# - Replace the monitored site/user with the new one (CARTAGO / SEDE INGETRONIK / YISUS)
# - Flip ESTADO from DESCONECTADO (red) to CONECTADO (green)
# - Drop the two extra rows that are no longer part of the report
# - Re-fit the CIUDAD/USUARIO/ESTADO columns to the new (shorter) content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2019-10-16 11:49:42"

# Title banner (merged A1:E1)
$ws.Range("A1").Value = "Reporte generado: " + $newTimestamp

# Remove the two stale monitoring rows (6 and 7); row 5 becomes the only data row
$ws.Rows("6:7").Delete()

# Update the remaining data row with the new monitoring entry
$ws.Range("A5").Value = "CARTAGO"
$ws.Range("B5").Value = "SEDE INGETRONIK"
$ws.Range("C5").Value = "YISUS"
$ws.Range("D5").Value = $newTimestamp
$ws.Range("E5").Value = "CONECTADO"

# CONECTADO -> green fill (was red for DESCONECTADO); 0x45A236 packed as BGR for COM
$ws.Range("E5").Interior.Color = 3580485

# Re-fit columns to the new content (SEDE / USUARIO / ESTADO got shorter or longer)
$ws.Columns.Item(2).ColumnWidth = 17.833333
$ws.Columns.Item(3).ColumnWidth = 8.5
$ws.Columns.Item(5).ColumnWidth = 10.833333

# Keep the selection / used range in sync with the shrunk table
[void]$ws.Range("A4:E5").Select()
